$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "List Icon"

$newSheet.Range("A1").Value = "icon"
$newSheet.Range("A2").Value = "trash"
$newSheet.Range("A3").Value = "edit"

$newSheet.Activate()
$newSheet.Range("A4").Select() | Out-Null
